# Add violent crime data for 2022-08-16 across all affected worksheets.
# Each worksheet has year columns B:I (2015-2022) in row 1, crime categories
# in column A (rows 2+), and a "Total" row. Column I (2022 year-to-date totals)
# increases by the count of incidents recorded on 2022-08-16 for each category.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 4446
$ws.Range("I3").Value = 4647
$ws.Range("I4").Value = 1072
$ws.Range("I5").Value = 429
$ws.Range("I6").Value = 5065
$ws.Range("I7").Value = 15659

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 59
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 162
$ws.Range("I5").Value = 15
$ws.Range("I7").Value = 506

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I3").Value = 106
$ws.Range("I7").Value = 293

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I6").Value = 201
$ws.Range("I7").Value = 616

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 42
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 115
$ws.Range("I3").Value = 102
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 354

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I7").Value = 502
$ws.Range("I8").Value = 954
$ws.Range("I10").Value = 110
$ws.Range("I11").Value = 239
$ws.Range("I13").Value = 26
$ws.Range("I15").Value = 179
$ws.Range("I19").Value = 441
$ws.Range("I20").Value = 375
$ws.Range("I23").Value = 146
$ws.Range("I27").Value = 144
$ws.Range("I33").Value = 725
$ws.Range("I34").Value = 74
$ws.Range("I37").Value = 506
$ws.Range("I41").Value = 69
$ws.Range("I42").Value = 531
$ws.Range("I43").Value = 126
$ws.Range("I44").Value = 113
$ws.Range("I47").Value = 105
$ws.Range("I48").Value = 216
$ws.Range("I52").Value = 332
$ws.Range("I53").Value = 162
$ws.Range("I54").Value = 345
$ws.Range("I60").Value = 78
$ws.Range("I63").Value = 61
$ws.Range("I64").Value = 137
$ws.Range("I65").Value = 354
$ws.Range("I67").Value = 616
$ws.Range("I72").Value = 57
$ws.Range("I73").Value = 134
$ws.Range("I76").Value = 236
$ws.Range("I77").Value = 93
$ws.Range("I79").Value = 434
$ws.Range("I81").Value = 15
$ws.Range("I83").Value = 318
$ws.Range("I84").Value = 129
$ws.Range("I85").Value = 700
$ws.Range("I88").Value = 145
$ws.Range("I94").Value = 146
$ws.Range("I95").Value = 260
$ws.Range("I96").Value = 163
$ws.Range("I99").Value = 293
$ws.Range("I101").Value = 15659

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 318

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 91
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 260

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 270
$ws.Range("I6").Value = 228
$ws.Range("I7").Value = 725

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 71
$ws.Range("I7").Value = 345

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 288
$ws.Range("I6").Value = 272

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 163
$ws.Range("I5").Value = 10
$ws.Range("I7").Value = 441

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 216

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I3").Value = 54
$ws.Range("I4").Value = 28
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 236

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 180
$ws.Range("I3").Value = 281
$ws.Range("I7").Value = 700

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 183
$ws.Range("I5").Value = 19
$ws.Range("I7").Value = 531

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("I5").Value = 9
$ws.Range("I6").Value = 26

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 41
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 137
$ws.Range("I6").Value = 125
$ws.Range("I7").Value = 434

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 42
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 105
$ws.Range("I7").Value = 375

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 332

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I4").Value = 20
$ws.Range("I7").Value = 239

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 145

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 299
$ws.Range("I3").Value = 272
$ws.Range("I7").Value = 954

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I2").Value = 41
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 144

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 25
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I4").Value = 9
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 162

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 167
$ws.Range("I3").Value = 157
$ws.Range("I6").Value = 128
$ws.Range("I7").Value = 502

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("I2").Value = 9
$ws.Range("I6").Value = 15
